$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 8 (ano = 2025) with the revised customer metrics
$ws.Range("C8").Value = 977
$ws.Range("E8").Value = 814
$ws.Range("G8").Value = 83.31627430910952
$ws.Range("H8").Value = 16.68372569089048
